# Insert two new rows before row 86 and push all the subsequent rows down
# (this mirrors a new weekly price report being added at the top of the
# "Pepino ensalada" block, with the older reports shifting down by two rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("86:87").Insert()

# New row 86: Primera quality, week of 2022-03-25 (serial 44651)
$ws.Cells.Item(86, 1).Value = 2
$ws.Cells.Item(86, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = 44651
$ws.Cells.Item(86, 5).Value = 4
$ws.Cells.Item(86, 6).Value = 100112043
$ws.Cells.Item(86, 7).Value = "Pepino ensalada"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 300
$ws.Cells.Item(86, 11).Value = 15000
$ws.Cells.Item(86, 12).Value = 16000
$ws.Cells.Item(86, 13).Value = 15500
$ws.Cells.Item(86, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(86, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(86, 16).Value = 221
$ws.Cells.Item(86, 17).Value = 70
$ws.Cells.Item(86, 18).Value = "Hortaliza"

# New row 87: Segunda quality, same week (serial 44651)
$ws.Cells.Item(87, 1).Value = 2
$ws.Cells.Item(87, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44651
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = 100112043
$ws.Cells.Item(87, 7).Value = "Pepino ensalada"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Segunda"
$ws.Cells.Item(87, 10).Value = 248
$ws.Cells.Item(87, 11).Value = 12000
$ws.Cells.Item(87, 12).Value = 13000
$ws.Cells.Item(87, 13).Value = 12516
$ws.Cells.Item(87, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(87, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(87, 16).Value = 125
$ws.Cells.Item(87, 17).Value = 100
$ws.Cells.Item(87, 18).Value = "Hortaliza"
